$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I21").Value = "sd"
$ws.Range("J21").Value = "Statement-non-opinion"
$ws.Range("I23").Value = "ba"
$ws.Range("J23").Value = "Appreciation"
$ws.Range("I29").Value = "%"
$ws.Range("J29").Value = "Uninterpretable"
$ws.Range("I37").Value = "ba"
$ws.Range("J37").Value = "Appreciation"
$ws.Range("I44").Value = "sd"
$ws.Range("J44").Value = "Statement-non-opinion"
$ws.Range("I54").Value = "aa"
$ws.Range("J54").Value = "Agree/Accept"
$ws.Range("I67").Value = "aa"
$ws.Range("J67").Value = "Agree/Accept"
$ws.Range("I70").Value = "sv"
$ws.Range("J70").Value = "Statement-opinion"
$ws.Range("I87").Value = "sd"
$ws.Range("J87").Value = "Statement-non-opinion"
$ws.Range("I89").Value = "sd"
$ws.Range("J89").Value = "Statement-non-opinion"
$ws.Range("I97").Value = "b"
$ws.Range("J97").Value = "Acknowledge (Backchannel)"
$ws.Range("I120").Value = "aa"
$ws.Range("J120").Value = "Agree/Accept"
$ws.Range("I131").Value = "aa"
$ws.Range("J131").Value = "Agree/Accept"
$ws.Range("I139").Value = "aa"
$ws.Range("J139").Value = "Agree/Accept"
$ws.Range("I144").Value = "aa"
$ws.Range("J144").Value = "Agree/Accept"
$ws.Range("I162").Value = "b"
$ws.Range("J162").Value = "Acknowledge (Backchannel)"
$ws.Range("I170").Value = "aa"
$ws.Range("J170").Value = "Agree/Accept"
$ws.Range("I173").Value = "sd"
$ws.Range("J173").Value = "Statement-non-opinion"
$ws.Range("I178").Value = "ba"
$ws.Range("J178").Value = "Appreciation"
$ws.Range("I183").Value = "ba"
$ws.Range("J183").Value = "Appreciation"
$ws.Range("I206").Value = "sv"
$ws.Range("J206").Value = "Statement-opinion"
$ws.Range("I210").Value = "aa"
$ws.Range("J210").Value = "Agree/Accept"
$ws.Range("I216").Value = "sd"
$ws.Range("J216").Value = "Statement-non-opinion"
$ws.Range("I227").Value = "ba"
$ws.Range("J227").Value = "Appreciation"
$ws.Range("I252").Value = "ba"
$ws.Range("J252").Value = "Appreciation"
$ws.Range("I254").Value = "sv"
$ws.Range("J254").Value = "Statement-opinion"
$ws.Range("I260").Value = "sd"
$ws.Range("J260").Value = "Statement-non-opinion"
$ws.Range("I272").Value = "sd"
$ws.Range("J272").Value = "Statement-non-opinion"
$ws.Range("I276").Value = "aa"
$ws.Range("J276").Value = "Agree/Accept"
$ws.Range("I284").Value = "aa"
$ws.Range("J284").Value = "Agree/Accept"
$ws.Range("I292").Value = "sd"
$ws.Range("J292").Value = "Statement-non-opinion"
$ws.Range("I296").Value = "aa"
$ws.Range("J296").Value = "Agree/Accept"
$ws.Range("I297").Value = "sd"
$ws.Range("J297").Value = "Statement-non-opinion"
$ws.Range("I306").Value = "sd"
$ws.Range("J306").Value = "Statement-non-opinion"
$ws.Range("I320").Value = "aa"
$ws.Range("J320").Value = "Agree/Accept"
$ws.Range("I338").Value = "b"
$ws.Range("J338").Value = "Acknowledge (Backchannel)"
$ws.Range("I350").Value = "%"
$ws.Range("J350").Value = "Uninterpretable"
$ws.Range("I353").Value = "sv"
$ws.Range("J353").Value = "Statement-opinion"
$ws.Range("I374").Value = "sd"
$ws.Range("J374").Value = "Statement-non-opinion"
$ws.Range("I387").Value = "ba"
$ws.Range("J387").Value = "Appreciation"
$ws.Range("I391").Value = "b"
$ws.Range("J391").Value = "Acknowledge (Backchannel)"
$ws.Range("I396").Value = "sv"
$ws.Range("J396").Value = "Statement-opinion"
$ws.Range("I410").Value = "ba"
$ws.Range("J410").Value = "Appreciation"
$ws.Range("I440").Value = "aa"
$ws.Range("J440").Value = "Agree/Accept"
$ws.Range("I441").Value = "sd"
$ws.Range("J441").Value = "Statement-non-opinion"
$ws.Range("I445").Value = "b"
$ws.Range("J445").Value = "Acknowledge (Backchannel)"
$ws.Range("I446").Value = "sd"
$ws.Range("J446").Value = "Statement-non-opinion"
$ws.Range("I448").Value = "ba"
$ws.Range("J448").Value = "Appreciation"
$ws.Range("I453").Value = "b"
$ws.Range("J453").Value = "Acknowledge (Backchannel)"
$ws.Range("I488").Value = "sd"
$ws.Range("J488").Value = "Statement-non-opinion"
$ws.Range("I490").Value = "b"
$ws.Range("J490").Value = "Acknowledge (Backchannel)"
$ws.Range("I494").Value = "sv"
$ws.Range("J494").Value = "Statement-opinion"
$ws.Range("I514").Value = "b"
$ws.Range("J514").Value = "Acknowledge (Backchannel)"
$ws.Range("I534").Value = "sv"
$ws.Range("J534").Value = "Statement-opinion"
$ws.Range("I548").Value = "sd"
$ws.Range("J548").Value = "Statement-non-opinion"
$ws.Range("I554").Value = "sd"
$ws.Range("J554").Value = "Statement-non-opinion"
$ws.Range("I566").Value = "b"
$ws.Range("J566").Value = "Acknowledge (Backchannel)"
$ws.Range("I575").Value = "ba"
$ws.Range("J575").Value = "Appreciation"
$ws.Range("I597").Value = "sd"
$ws.Range("J597").Value = "Statement-non-opinion"
$ws.Range("I599").Value = "sv"
$ws.Range("J599").Value = "Statement-opinion"
$ws.Range("I601").Value = "sv"
$ws.Range("J601").Value = "Statement-opinion"
$ws.Range("I608").Value = "sv"
$ws.Range("J608").Value = "Statement-opinion"
$ws.Range("I612").Value = "sv"
$ws.Range("J612").Value = "Statement-opinion"
$ws.Range("I613").Value = "%"
$ws.Range("J613").Value = "Uninterpretable"
$ws.Range("I615").Value = "sv"
$ws.Range("J615").Value = "Statement-opinion"
$ws.Range("I619").Value = "aa"
$ws.Range("J619").Value = "Agree/Accept"
$ws.Range("I639").Value = "b"
$ws.Range("J639").Value = "Acknowledge (Backchannel)"
$ws.Range("I641").Value = "sv"
$ws.Range("J641").Value = "Statement-opinion"
$ws.Range("I649").Value = "aa"
$ws.Range("J649").Value = "Agree/Accept"
$ws.Range("I656").Value = "sv"
$ws.Range("J656").Value = "Statement-opinion"
$ws.Range("I667").Value = "sd"
$ws.Range("J667").Value = "Statement-non-opinion"
